{"js": "// Merge the three template-placeholder runs for the \"input date\" table cell\n// into a single run that applies Jinja's `date` filter (with a leading space\n// before the optional paper-date parenthetical).\nconst oldText =\n  \"{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}\" +\n  \"({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}\";\nconst newText =\n  '{% if inputDateHeader %}{{ inputDateHeader | date(\"dd.MM.YYYY\") }}' +\n  '{% if paperInputDateHeader %} ({{ paperInputDateHeader | date(\"dd.MM.YYYY\") }})' +\n  \"{% else %}{% endif %}{% else %}-{% endif %}\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target placeholder text to update.\");\n}\n\nresults.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Merge the three template-placeholder runs for the \"input date\" table cell\n# into a single run that applies Jinja's `date` filter (with a leading space\n# before the optional paper-date parenthetical).\n$d = $word.ActiveDocument\n\n$oldText = '{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}'\n$newText = '{% if inputDateHeader %}{{ inputDateHeader | date(\"dd.MM.YYYY\") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date(\"dd.MM.YYYY\") }}){% else %}{% endif %}{% else %}-{% endif %}'\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target placeholder text to update.\"\n}\n\n# Assign directly to Range.Text (rather than Find.Replacement.Text) so Word's\n# smart-quote AutoCorrect does not mangle the literal straight quotes in the\n# Jinja filter arguments. This also collapses the matched range's three runs\n# into a single run that inherits the first run's formatting.\n$range.Text = $newText\n"}
